# Apply Sun Jun 18 2023 cryptos-list refresh: update Price/Volume columns for
# existing rows (2-33) and shift the coin rows 34-51 (Frax drops off the
# tracked list, every following coin moves up one slot, Aave is newly added
# at the bottom).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-33: update Price (D) and Volume(1h) (E) columns only ---
$deUpdates = @(
    @{ Row=2; D="26.502.16"; E="  -0.38%  " },
    @{ Row=3; D="1.728.58"; E="  -0.93%  " },
    @{ Row=4; D=$null; E="  +0.17%  " },
    @{ Row=5; D="246.27"; E="  -0.39%  " },
    @{ Row=6; D=$null; E="  +0.09%  " },
    @{ Row=7; D="0.4819"; E="  +0.09%  " },
    @{ Row=8; D="0.2661"; E="  -1.40%  " },
    @{ Row=9; D="0.06216"; E="  -0.81%  " },
    @{ Row=10; D="1.728.35"; E="  -0.90%  " },
    @{ Row=11; D="0.07073"; E="  -0.65%  " },
    @{ Row=12; D="15.56"; E="  -1.88%  " },
    @{ Row=13; D="4.583"; E="  +1.49%  " },
    @{ Row=14; D="0.6086"; E="  -2.47%  " },
    @{ Row=15; D="77.28"; E="  -0.38%  " },
    @{ Row=16; D="1.001"; E="  +0.06%  " },
    @{ Row=17; D="26.500.39"; E="  -0.34%  " },
    @{ Row=18; D=$null; E="  +0.16%  " },
    @{ Row=19; D="0.000007190"; E="  +4.02%  " },
    @{ Row=20; D=$null; E="  -1.73%  " },
    @{ Row=21; D="1.952.42"; E="  -0.73%  " },
    @{ Row=22; D="4.498"; E="  -3.01%  " },
    @{ Row=23; D="8.760"; E="  -0.87%  " },
    @{ Row=24; D="5.237"; E="  -2.33%  " },
    @{ Row=25; D="137.27"; E="  +0.97%  " },
    @{ Row=26; D="15.42"; E="  +0.09%  " },
    @{ Row=27; D="1.774"; E="  -2.62%  " },
    @{ Row=28; D="1.407"; E="  -2.04%  " },
    @{ Row=29; D="108.08"; E="  +0.95%  " },
    @{ Row=30; D="3.970"; E="  -1.07%  " },
    @{ Row=31; D="0.07987"; E="  +1.25%  " },
    @{ Row=32; D="3.685"; E="  -1.94%  " },
    @{ Row=33; D=$null; E="  -0.79%  " },
)

foreach ($item in $deUpdates) {
    if ($item.D -ne $null) {
        $dCell = $ws.Cells.Item($item.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $item.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# --- Rows 34-51: coin list shifted up one (Frax removed), new Aave row appended ---
$rowUpdates = @(
    @{ Row=34; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.620"; E="  +0.06%  " },
    @{ Row=35; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.002"; E="  +0.21%  " },
    @{ Row=36; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.6337"; E="  -1.89%  " },
    @{ Row=37; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.8892"; E="  -6.24%  " },
    @{ Row=38; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="2.013"; E="  +0.53%  " },
    @{ Row=39; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.394"; E="  -1.61%  " },
    @{ Row=40; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="1.002"; E="  -0.07%  " },
    @{ Row=41; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01501"; E="  -0.76%  " },
    @{ Row=42; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="101.83"; E="  -10.30%  " },
    @{ Row=43; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="5.500"; E="  -4.69%  " },
    @{ Row=44; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.3879"; E="  -1.22%  " },
    @{ Row=45; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="6.982"; E="  +3.71%  " },
    @{ Row=46; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="0.1183"; E="  -2.51%  " },
    @{ Row=47; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.05385"; E="  +1.04%  " },
    @{ Row=48; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="7.912"; E="  -1.08%  " },
    @{ Row=49; B="Elrond"; C="https://coinranking.com/coin/omwkOTglq+elrond-egld"; D="30.59"; E="  -0.60%  " },
    @{ Row=50; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="1.251"; E="  -1.64%  " },
    @{ Row=51; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="51.49"; E="  -0.41%  " },
)

foreach ($item in $rowUpdates) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.Style = "Normal"
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}
